$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep headers and column A (dates) as-is; replace the predicted values in
# column B with the refreshed model output.
$ws.Range("B2").Value = 174519.3831821379
$ws.Range("B3").Value = 172273.6505372123
$ws.Range("B4").Value = 179065.5350782446
$ws.Range("B5").Value = 184738.544669773
$ws.Range("B6").Value = 183341.1269387082
$ws.Range("B7").Value = 184311.2956594777
$ws.Range("B8").Value = 184886.9138254827

$wb.Save()
